# Progress Report - Week11: replace the trailing tab at the end of the
# "Week 11" paragraph with a new paragraph containing the financial update.
#
# The target OOXML turns the old:
#     ... scanning of barcodes.<space><tab><bookmarkStart/><bookmarkEnd/>
# into:
#     ... scanning of barcodes.<space></w:p>
#     <w:p>
#       <four runs of new text>
#       <bookmarkStart/><bookmarkEnd/>
#     </w:p>
#
# i.e. the tab becomes a paragraph break, the new sentences land in the
# freshly created (now last) paragraph as four separate runs, and the
# "_GoBack" bookmark (zero-width, around the former tab position) ends up
# right after the new text instead of right before it.

$d = $word.ActiveDocument

# 1) Turn the trailing tab character into a paragraph mark. This splits the
#    paragraph in two: the old text keeps its trailing space, and a brand
#    new (empty) last paragraph is created, carrying the "_GoBack" bookmark
#    at its start.
$find = $d.Content.Find
$find.Text = "^t"
$find.Replacement.Text = "^p"
$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null

# 2) Type the new sentences into that new, still-empty paragraph. Using
#    Selection.TypeText (rather than Range.InsertAfter/InsertBefore) keeps
#    each chunk as its own run instead of silently coalescing them into one
#    run, matching the four separate <w:r> elements in the diff.
$newParaStart = $d.Paragraphs($d.Paragraphs.Count).Range.Start

$chunks = @(
    "There has been an financial update on our project, as I said before I had a problem with my webcam which resulted in buying new webcam. This webcam cost me ap",
    "proximately 15 dollars. With ",
    "exception of the webcam have no not spent ",
    "any money since the last report."
)

$sel = $word.Selection
$pos = $newParaStart
foreach ($chunk in $chunks) {
    $sel.SetRange($pos, $pos)
    $sel.TypeText($chunk)
    $pos = $pos + $chunk.Length
}
$newTextEnd = $pos

# 3) Re-seat the "_GoBack" bookmark after the text we just typed. Bookmarks
#    placed with Start=End sitting right at the very end of the document
#    tend to mis-resolve, so pad the story out first, move the bookmark,
#    then remove the padding again.
$padding = "ZzPaddingPlaceholderZz"
$tailRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$tailRange.InsertAfter($padding)

$d.Bookmarks("_GoBack").Delete()
$bookmarkRange = $d.Range($newTextEnd, $newTextEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$padRange = $d.Range($newTextEnd, $newTextEnd + $padding.Length)
$padRange.Delete()
